# Update "想去人数" (interest count) values in column F on the "展览" and
# "全部类型" worksheets to match the newly scraped totals.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 1315
$ws.Range("F3").Value  = 1184
$ws.Range("F4").Value  = 14433
$ws.Range("F5").Value  = 16943
$ws.Range("F7").Value  = 129
$ws.Range("F8").Value  = 41
$ws.Range("F16").Value = 39
$ws.Range("F17").Value = 117
$ws.Range("F19").Value = 1292
$ws.Range("F22").Value = 54
$ws.Range("F23").Value = 35
$ws.Range("F25").Value = 6962
$ws.Range("F27").Value = 28
$ws.Range("F28").Value = 1147
$ws.Range("F29").Value = 17
$ws.Range("F31").Value = 5800
$ws.Range("F32").Value = 124
$ws.Range("F34").Value = 207
$ws.Range("F35").Value = 4937
$ws.Range("F36").Value = 27

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 1315
$ws.Range("F3").Value  = 1184
$ws.Range("F4").Value  = 14433
$ws.Range("F5").Value  = 16943
$ws.Range("F7").Value  = 129
$ws.Range("F8").Value  = 41
$ws.Range("F16").Value = 39
$ws.Range("F17").Value = 117
$ws.Range("F19").Value = 1292
$ws.Range("F23").Value = 54
$ws.Range("F24").Value = 35
$ws.Range("F26").Value = 6962
$ws.Range("F28").Value = 28
$ws.Range("F29").Value = 1147
$ws.Range("F30").Value = 17
$ws.Range("F33").Value = 5800
$ws.Range("F34").Value = 124
$ws.Range("F36").Value = 207
$ws.Range("F37").Value = 4937
$ws.Range("F38").Value = 27
